$d = $word.ActiveDocument
$insertRange = $d.Range(0, 0)
$newParasXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>「地方に特化 × Vibes &amp; Mood」</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>でやっていきたい</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>04/29</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>現状</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>インデックスページだけが、</w:t></w:r><w:r><w:t>explore Japan’s countryside</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>をアピールしている。その他のページは、単に</w:t></w:r><w:r><w:t>countryside</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>にある観光地を紹介している</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>目標</w:t></w:r><w:r><w:t xml:space="preserve">: Japan countryside, Japanese countryside </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>と言う検索で</w:t></w:r><w:r><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ページ目に入りたい</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>自分が考えている問題点</w:t></w:r><w:r><w:t>: SEO</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的に</w:t></w:r><w:r><w:t>countryside</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>感が弱い</w:t></w:r><w:r><w:t>. C</w:t></w:r><w:r><w:t>ountryside</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>の観光地を紹介しているだけでは弱いのかなと思っています</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>手段</w:t></w:r><w:r><w:t xml:space="preserve">: countryside, rural </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>等の言葉をたくさん使った記事を執筆すること</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>箱根</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>江ノ島</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>鎌倉</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>富士山付近</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>の5つに関する記事</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>04/30</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Best countryside in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>japan</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>みたいな感じで書いていく</w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ここで負けちゃあかんな、たくさん良い記事書いて、</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>楽しんでもらうのが目標</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>なん</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>だからここで諦めちゃいかんよ</w:t></w:r></w:p><w:p><w:r><w:t>5/1</w:t></w:r></w:p><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Countside_guide</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>下に、</w:t></w:r><w:r><w:t>articles</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>ディレクトリを入れて、</w:t></w:r><w:r><w:t>SEO</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>的にうまく</w:t></w:r><w:r><w:t>countryside</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>を認識してもらおう</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">5/2 </w:t></w:r></w:p><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>まあなんかよくわからんけど、とりあえず熱海周りの観光地を追加しよう</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">one by one </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>でいこう</w:t></w:r></w:p><w:p><w:r><w:t>Article</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>を書いていこう</w:t></w:r></w:p>'
$insertRange.InsertXML($newParasXml)
Write-Output $d.Paragraphs.Count
